$wb = $excel.ActiveWorkbook

$changes = @{
    "展览"     = @{ 2=386; 3=673; 4=112; 5=2107; 7=10963; 8=186; 9=164; 10=297; 12=10796; 13=432; 15=21; 16=752; 17=5407; 18=79; 19=3397 }
    "演出"     = @{ 3=563 }
    "全部类型" = @{ 2=386; 3=673; 5=112; 6=2107; 7=563; 10=10963; 11=186; 12=164; 13=297; 15=10796; 16=432; 18=21; 19=752; 20=5407; 21=79; 22=3397 }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $changes[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $value = $rowsMap[$row]
        $ws.Cells.Item($row, 6).Value = $value
    }
}
